$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)       # "总计"
$wsQ4ref = $wb.Worksheets.Item(2)   # positional anchor: existing "2021-Q4" sheet

# ---------------------------------------------------------------
# 1) Insert a brand-new "2022-Q3" worksheet right before "2021-Q4"
# ---------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Add($wsQ4ref)
$wsQ3.Name = "2022-Q3"

# NOTE: after Add(), the $wsQ4ref handle now tracks the freshly inserted
# sheet (position-bound), not the original "2021-Q4" sheet any more — so
# the real "2021-Q4" worksheet must be re-acquired by name.
$wsQ4 = $wb.Worksheets.Item("2021-Q4")

# Match existing look & feel: copy the bold/bordered header style (s=2)
# from the "总计" sheet's header row, and the matching style used for
# column A, onto the new sheet.
$ws1.Range("B1:D1").Copy() | Out-Null
$wsQ3.Range("B1:H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws1.Range("A2").Copy() | Out-Null
$wsQ3.Range("A2:A3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Header row
$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"

# Helper values are written as TEXT (they must stay text, not numbers, to
# match the source data) by briefly switching the cell to the "@" (Text)
# number format, assigning the literal string, then restoring the cell's
# style to Normal so no stray formatting is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - 南方金融主题灵活配置混合A
$wsQ3.Range("A2").Value = 0
Set-TextValue $wsQ3.Range("B2") "004702"
$wsQ3.Range("C2").Value = "南方金融主题灵活配置混合A"
Set-TextValue $wsQ3.Range("D2") "11.74"
Set-TextValue $wsQ3.Range("E2") "92.15"
Set-TextValue $wsQ3.Range("F2") "3.44"
Set-TextValue $wsQ3.Range("G2") "0.4039"
$wsQ3.Range("H2").Value = 7

# Row 3 - 南方金融主题灵活配置混合C
$wsQ3.Range("A3").Value = 1
Set-TextValue $wsQ3.Range("B3") "013500"
$wsQ3.Range("C3").Value = "南方金融主题灵活配置混合C"
Set-TextValue $wsQ3.Range("D3") "7.39"
Set-TextValue $wsQ3.Range("E3") "92.15"
Set-TextValue $wsQ3.Range("F3") "3.44"
Set-TextValue $wsQ3.Range("G3") "0.2542"
$wsQ3.Range("H3").Value = 7

# ---------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert a new row for 2022-Q3
#    above the existing 2021-Q4 row (which shifts down to row 3
#    with its A-value bumped from 0 to 1).
# ---------------------------------------------------------------
$b2old = $ws1.Range("B2").Value2
$c2old = $ws1.Range("C2").Value2
$d2old = $ws1.Range("D2").Value2

# Duplicate A2's format onto A3 so the shifted row keeps style s=2
$ws1.Range("A2").Copy() | Out-Null
$ws1.Range("A3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = $b2old
$ws1.Range("C3").Value = $c2old
$ws1.Range("D3").Value = $d2old

$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 2
$ws1.Range("D2").Value = 0.66

# ---------------------------------------------------------------
# 3) Keep the originally-selected tab ("2021-Q4") selected, since it
#    carried over from the source workbook and this edit doesn't
#    touch the view state otherwise.
# ---------------------------------------------------------------
$wsQ4.Activate()

Write-Output "edit complete"
